$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2248.875
$ws.Range("I132").Value = 1141.6428
$ws.Range("K132").Value = 3424.9284
$ws.Range("M132").Value = -894.9284000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2707.1177
$ws.Range("J137").Value = 4142.2856
$ws.Range("L137").Value = 12426.8568
$ws.Range("N137").Value = -17526.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2251.375
$ws.Range("I2").Value = 2251.375
$ws.Range("K2").Value = 2251.375
$ws.Range("M2").Value = -2138.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 16770
$ws.Range("I36").Value = 16770
$ws.Range("K36").Value = 16770
$ws.Range("M36").Value = -16424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6867.4443
$ws.Range("I74").Value = 7538.375
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 7538.375
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -6664.375
$ws.Range("N74").Value = -3248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6867.4443
$ws.Range("I77").Value = 7538.375
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 37691.875
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -33323.875
$ws.Range("N77").Value = -16236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2087.5789
$ws.Range("I102").Value = 2087.5789
$ws.Range("K102").Value = 2087.5789
$ws.Range("M102").Value = -465.5789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2251.375
$ws.Range("I116").Value = 2251.375
$ws.Range("K116").Value = 2251.375
$ws.Range("M116").Value = 42.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3805.8235
$ws.Range("I132").Value = 3558.4167
$ws.Range("K132").Value = 10675.2501
$ws.Range("M132").Value = -8145.250100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2251.375
$ws.Range("I3").Value = 2251.375
$ws.Range("K3").Value = 2251.375
$ws.Range("M3").Value = -2137.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4460.4287
$ws.Range("I94").Value = 2531.4546
$ws.Range("J94").Value = 11533.333
$ws.Range("K94").Value = 2531.4546
$ws.Range("L94").Value = 11533.333
$ws.Range("M94").Value = -2080.4546
$ws.Range("N94").Value = -12435.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1552.75
$ws.Range("I107").Value = 1155.5
$ws.Range("J107").Value = 1950
$ws.Range("K107").Value = 1155.5
$ws.Range("L107").Value = 1950
$ws.Range("M107").Value = 764.5
$ws.Range("N107").Value = -5790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3326.6667
$ws.Range("I134").Value = 3326.6667
$ws.Range("K134").Value = 9980.000100000001
$ws.Range("M134").Value = -7445.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 73474.71000000001
$ws.Range("I16").Value = 101494.7
$ws.Range("K16").Value = 101494.7
$ws.Range("M16").Value = -101207.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2512.6667
$ws.Range("J31").Value = 3088.2
$ws.Range("L31").Value = 3088.2
$ws.Range("N31").Value = -3678.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2512.6667
$ws.Range("J34").Value = 3088.2
$ws.Range("L34").Value = 3088.2
$ws.Range("N34").Value = -3492.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1729.4117
$ws.Range("I58").Value = 1800.7778
$ws.Range("K58").Value = 1800.7778
$ws.Range("M58").Value = -1597.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 73474.71000000001
$ws.Range("I113").Value = 101494.7
$ws.Range("K113").Value = 101494.7
$ws.Range("M113").Value = -99324.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 79422.69500000001
$ws.Range("I132").Value = 112944.11
$ws.Range("K132").Value = 338832.33
$ws.Range("M132").Value = -336302.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4041.0833
$ws.Range("I134").Value = 4249.5
$ws.Range("K134").Value = 12748.5
$ws.Range("M134").Value = -10213.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1729.4117
$ws.Range("I136").Value = 1800.7778
$ws.Range("K136").Value = 5402.3334
$ws.Range("M136").Value = -2852.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 716.3333
$ws.Range("J121").Value = 716.3333
$ws.Range("L121").Value = 2148.9999
$ws.Range("N121").Value = -4768.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 358.16666
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 566.3333
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 566.3333
$ws.Range("M2").Value = -37
$ws.Range("N2").Value = -792.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("N70").Value = -5540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5000
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("N73").Value = -6872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2889
$ws.Range("I97").Value = 2128.7144
$ws.Range("J97").Value = 5550
$ws.Range("K97").Value = 2128.7144
$ws.Range("L97").Value = 5550
$ws.Range("M97").Value = -1632.7144
$ws.Range("N97").Value = -6542

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3173.4
$ws.Range("I122").Value = 3241.75
$ws.Range("K122").Value = 9725.25
$ws.Range("M122").Value = -7275.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2043.05
$ws.Range("I132").Value = 1580.1765
$ws.Range("K132").Value = 4740.529500000001
$ws.Range("M132").Value = -2210.529500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 35326
$ws.Range("J136").Value = 35326
$ws.Range("L136").Value = 105978
$ws.Range("N136").Value = -111078

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 7559
$ws.Range("I82").Value = 2588.5
$ws.Range("K82").Value = 2588.5
$ws.Range("M82").Value = -2227.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 7559
$ws.Range("I85").Value = 2588.5
$ws.Range("K85").Value = 2588.5
$ws.Range("M85").Value = -1340.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2810
$ws.Range("I100").Value = 2736
$ws.Range("K100").Value = 2736
$ws.Range("M100").Value = -2195

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3990.4285
$ws.Range("I132").Value = 2483.75
$ws.Range("K132").Value = 7451.25
$ws.Range("M132").Value = -4921.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 906.75
$ws.Range("I100").Value = 252.8
$ws.Range("J100").Value = 1996.6666
$ws.Range("K100").Value = 505.6
$ws.Range("L100").Value = 3993.3332
$ws.Range("M100").Value = 35.39999999999998
$ws.Range("N100").Value = -5075.3332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2050.5881
$ws.Range("I132").Value = 1155.1666
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 3465.4998
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -935.4998000000001
$ws.Range("N132").Value = -17658.8
